$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 267
$ws1.Range("F4").Value = 70
$ws1.Range("F5").Value = 248
$ws1.Range("F7").Value = 77
$ws1.Range("F12").Value = 104
$ws1.Range("F13").Value = 2340
$ws1.Range("F16").Value = 520
$ws1.Range("F17").Value = 536
$ws1.Range("F22").Value = 1812
$ws1.Range("F23").Value = 3939
$ws1.Range("F24").Value = 29
$ws1.Range("F26").Value = 1172
$ws1.Range("F28").Value = 2080
$ws1.Range("F30").Value = 466
$ws1.Range("F32").Value = 103
$ws1.Range("F36").Value = 685
$ws1.Range("F38").Value = 406

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 29

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 267
$ws4.Range("F4").Value = 70
$ws4.Range("F5").Value = 248
$ws4.Range("F7").Value = 77
$ws4.Range("F12").Value = 104
$ws4.Range("F13").Value = 2340
$ws4.Range("F15").Value = 29
$ws4.Range("F17").Value = 520
$ws4.Range("F18").Value = 536
$ws4.Range("F23").Value = 1812
$ws4.Range("F24").Value = 3939
$ws4.Range("F25").Value = 29
$ws4.Range("F27").Value = 1172
$ws4.Range("F29").Value = 2080
$ws4.Range("F31").Value = 466
$ws4.Range("F33").Value = 103
$ws4.Range("F37").Value = 685
$ws4.Range("F39").Value = 406

